# Auto-generated edit script: updates crypto price/volume table cells
# to match the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.135.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "'3.274.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'587.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("D6").Value = "'186.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.53%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("E9").Value = "  +4.35%  "
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("D11").Value = "'0.418"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "'3.838.95"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").Value = "'28.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.48%  "
$ws.Range("D15").Value = "'68.114.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("D17").Value = "'3.268.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "'5.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "'13.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("D20").Value = "'381.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.76%  "
$ws.Range("D21").Value = "'7.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "'71.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").Value = "'0.515"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("D26").Value = "'9.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("D27").Value = "'0.188"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.98%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("E29").Value = "  +2.51%  "
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").Value = "'22.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("D32").Value = "'7.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.22%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("E35").Value = "  +2.50%  "
$ws.Range("D36").Value = "'162.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.50%  "
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("E38").Value = "  -2.20%  "
$ws.Range("D39").Value = "'6.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'26.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("D41").Value = "'4.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.51%  "
$ws.Range("D42").Value = "'2.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "'0.0690"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.14%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'41.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.90%  "
$ws.Range("D45").Value = "'25.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "'346.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "'2.649.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.30%  "
$ws.Range("D48").Value = "'0.0285"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("E49").Value = "  +4.20%  "
$ws.Range("D50").Value = "'1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("D51").Value = "'0.103"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.25%  "
